$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.667.51'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '1.839.25'
$ws.Range("E3").Value = '  +1.63%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.24'
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.610'
$ws.Range("E6").Value = '  +1.96%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.26'
$ws.Range("E8").Value = '  +16.70%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.301'
$ws.Range("E9").Value = '  +3.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0684'
$ws.Range("E10").Value = '  +0.49%  '
$ws.Range("E11").Value = '  +3.57%  '
$ws.Range("D12").Value = '2.104.80'
$ws.Range("E12").Value = '  +1.64%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.837.31'
$ws.Range("E13").Value = '  +1.69%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.26'
$ws.Range("E14").Value = '  -0.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.67'
$ws.Range("E15").Value = '  +5.71%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.650'
$ws.Range("E16").Value = '  +3.20%  '
$ws.Range("D17").Value = '34.664.36'
$ws.Range("E17").Value = '  +0.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.49'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.75'
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("D20").Value = '0.0₃0783'
$ws.Range("E20").Value = '  +1.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.01'
$ws.Range("E21").Value = '  +7.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.71'
$ws.Range("E22").Value = '  +14.88%  '
$ws.Range("E23").Value = '  +0.26%  '
$ws.Range("E24").Value = '  -1.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.91'
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.84'
$ws.Range("E26").Value = '  -0.41%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.72'
$ws.Range("E27").Value = '  +2.74%  '
$ws.Range("E28").Value = '  +0.19%  '
$ws.Range("E29").Value = '  +0.39%  '
$ws.Range("E30").Value = '  +5.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.90'
$ws.Range("E31").Value = '  +2.06%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.97'
$ws.Range("E32").Value = '  +1.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0526'
$ws.Range("E33").Value = '  +1.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.86'
$ws.Range("E34").Value = '  +3.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '89.69'
$ws.Range("E35").Value = '  +10.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.660'
$ws.Range("E36").Value = '  +1.13%  '
$ws.Range("D37").Value = '1.337.97'
$ws.Range("E37").Value = '  -1.86%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.42'
$ws.Range("E38").Value = '  +2.83%  '
$ws.Range("E39").Value = '  +0.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0192'
$ws.Range("E40").Value = '  +3.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.91'
$ws.Range("E41").Value = '  +11.33%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.983'
$ws.Range("E42").Value = '  +5.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.24'
$ws.Range("E43").Value = '  +7.36%  '
$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.44'
$ws.Range("E44").Value = '  +0.68%  '
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.82'
$ws.Range("E45").Value = '  +1.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0520'
$ws.Range("E46").Value = '  +4.10%  '
$ws.Range("D47").Value = '2.004.75'
$ws.Range("E47").Value = '  +1.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.98'
$ws.Range("E48").Value = '  +2.72%  '
$ws.Range("E49").Value = '  +0.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.04'
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0611'
$ws.Range("E51").Value = '  +0.89%  '
